# Apply cryptos list update (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.741.68'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.627.90'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.81%  '
$ws.Range('D5').Value = "'214.08"
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').Value = "'0.500"
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = "'19.58"
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.853.40'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.630.43'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').Value = "'0.551"
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = "'62.69"
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '25.746.85'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = "'0.998"
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = "'4.43"
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').Value = "'190.81"
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = "'6.26"
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').Value = "'0.998"
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').Value = "'1.81"
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('D26').Value = "'142.02"
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('D28').Value = "'6.81"
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').Value = "'15.46"
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').Value = "'1.23"
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = "'0.0493"
$ws.Range('E31').Value = '  +2.19%  '
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = "'3.22"
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('D34').Value = "'1.58"
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').Value = "'0.901"
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('D37').Value = '1.144.53'
$ws.Range('E37').Value = '  +3.88%  '
$ws.Range('D38').Value = "'0.542"
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').Value = "'0.997"
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('D43').Value = "'5.60"
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('D45').Value = "'0.801"
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').Value = '1.762.98'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'55.23"
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.46"
$ws.Range('E48').Value = '  +7.50%  '
$ws.Range('D49').Value = "'0.0512"
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = "'0.416"
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'2.32"
$ws.Range('E51').Value = '  -1.59%  '
